# Insert a new "Position (Enum)" sub-bullet listing the enum values,
# right after the existing "Position (Enum)" bullet and before "GameScore".
$d = $word.ActiveDocument

$rng = $d.Content
$found = $rng.Find.Execute("Position (Enum)", $true, $false, $false, $false, $false,
                            $true, 1, $false, "", 0)

if ($found) {
    # Collapse to the end of the found text (end of the "Position (Enum)" run,
    # just before its paragraph mark) and insert a brand-new paragraph there.
    # Word clones the pPr (style/numbering) of the paragraph it splits from,
    # which gives us the correct ListParagraph / ilvl=1 / numId=1 formatting.
    $rng.Collapse(0)
    $rng.InsertParagraphAfter()

    # Locate the freshly inserted (still-empty) paragraph: it is the one
    # immediately following the "Position (Enum)" paragraph.
    $paragraphs = $d.Paragraphs
    $count = $paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $para = $paragraphs.Item($i)
        if ($para.Range.Text -eq "Position (Enum)`r") {
            $newPara = $paragraphs.Item($i + 1)
            $newPara.Range.Text = "Pitcher, Catcher, FirstBase, SecondBase,ShortStop, ThirdBase, LeftField, RightField, CenterField"
            break
        }
    }
}
